$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44789
$ws.Range("J2").Value2 = 90
$ws.Range("K2").Value2 = 24000
$ws.Range("L2").Value2 = 24000
$ws.Range("M2").Value2 = 24000
$ws.Range("P2").Value2 = 1600
$ws.Range("D3").Value2 = 44407
$ws.Range("J3").Value2 = 90
$ws.Range("K3").Value2 = 25000
$ws.Range("L3").Value2 = 25000
$ws.Range("M3").Value2 = 25000
$ws.Range("P3").Value2 = 1667
$ws.Range("D4").Value2 = 44838
$ws.Range("J4").Value2 = 80
$ws.Range("K4").Value2 = 22000
$ws.Range("L4").Value2 = 22000
$ws.Range("M4").Value2 = 22000
$ws.Range("P4").Value2 = 1467
$ws.Range("D5").Value2 = 44810
$ws.Range("J5").Value2 = 110
$ws.Range("K5").Value2 = 22000
$ws.Range("L5").Value2 = 22000
$ws.Range("M5").Value2 = 22000
$ws.Range("P5").Value2 = 1467
$ws.Range("D6").Value2 = 44781
$ws.Range("J6").Value2 = 70
$ws.Range("K6").Value2 = 24000
$ws.Range("L6").Value2 = 24000
$ws.Range("M6").Value2 = 24000
$ws.Range("P6").Value2 = 1600
$ws.Range("D7").Value2 = 44817
$ws.Range("J7").Value2 = 90
$ws.Range("K7").Value2 = 23000
$ws.Range("L7").Value2 = 23000
$ws.Range("M7").Value2 = 23000
$ws.Range("P7").Value2 = 1533
$ws.Range("D8").Value2 = 44365
$ws.Range("J8").Value2 = 80
$ws.Range("K8").Value2 = 25000
$ws.Range("M8").Value2 = 25000
$ws.Range("P8").Value2 = 1667
$ws.Range("D9").Value2 = 44778
$ws.Range("J9").Value2 = 120
$ws.Range("K9").Value2 = 24000
$ws.Range("L9").Value2 = 24000
$ws.Range("M9").Value2 = 24000
$ws.Range("P9").Value2 = 1600
$ws.Range("D10").Value2 = 44764
$ws.Range("J10").Value2 = 90
$ws.Range("K10").Value2 = 24000
$ws.Range("L10").Value2 = 24000
$ws.Range("M10").Value2 = 24000
$ws.Range("P10").Value2 = 1600
$ws.Range("D11").Value2 = 44827
$ws.Range("J11").Value2 = 90
$ws.Range("K11").Value2 = 22000
$ws.Range("L11").Value2 = 22000
$ws.Range("M11").Value2 = 22000
$ws.Range("P11").Value2 = 1467
$ws.Range("D12").Value2 = 44771
$ws.Range("D14").Value2 = 44792
$ws.Range("J14").Value2 = 120
$ws.Range("K14").Value2 = 24000
$ws.Range("L14").Value2 = 24000
$ws.Range("M14").Value2 = 24000
$ws.Range("P14").Value2 = 1600
$ws.Range("D15").Value2 = 44754
$ws.Range("K15").Value2 = 25000
$ws.Range("L15").Value2 = 25000
$ws.Range("M15").Value2 = 25000
$ws.Range("P15").Value2 = 1667
$ws.Range("D18").Value2 = 44750
$ws.Range("J18").Value2 = 90
$ws.Range("K18").Value2 = 25000
$ws.Range("L18").Value2 = 25000
$ws.Range("M18").Value2 = 25000
$ws.Range("P18").Value2 = 1667
$ws.Range("D19").Value2 = 44418
$ws.Range("J19").Value2 = 90
$ws.Range("K19").Value2 = 25000
$ws.Range("L19").Value2 = 25000
$ws.Range("M19").Value2 = 25000
$ws.Range("P19").Value2 = 1667
$ws.Range("D20").Value2 = 44799
$ws.Range("J20").Value2 = 80
$ws.Range("D21").Value2 = 44831
$ws.Range("K21").Value2 = 25000
$ws.Range("L21").Value2 = 25000
$ws.Range("M21").Value2 = 25000
$ws.Range("P21").Value2 = 1667
$ws.Range("D22").Value2 = 44761
$ws.Range("J22").Value2 = 100
$ws.Range("K22").Value2 = 23000
$ws.Range("M22").Value2 = 24000
$ws.Range("P22").Value2 = 1600
$ws.Range("D23").Value2 = 44740
$ws.Range("J23").Value2 = 90
$ws.Range("K23").Value2 = 25000
$ws.Range("L23").Value2 = 25000
$ws.Range("M23").Value2 = 25000
$ws.Range("P23").Value2 = 1667
$ws.Range("D24").Value2 = 44400
$ws.Range("J24").Value2 = 80
$ws.Range("K24").Value2 = 25000
$ws.Range("L24").Value2 = 25000
$ws.Range("M24").Value2 = 25000
$ws.Range("P24").Value2 = 1667
$ws.Range("D25").Value2 = 44819
$ws.Range("J25").Value2 = 70
$ws.Range("K25").Value2 = 22000
$ws.Range("L25").Value2 = 22000
$ws.Range("M25").Value2 = 22000
$ws.Range("P25").Value2 = 1467
$ws.Range("D26").Value2 = 44757
$ws.Range("J26").Value2 = 80
$ws.Range("D27").Value2 = 44806
$ws.Range("J27").Value2 = 70
$ws.Range("K27").Value2 = 23000
$ws.Range("L27").Value2 = 23000
$ws.Range("M27").Value2 = 23000
$ws.Range("P27").Value2 = 1533
